$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy style from existing header cell (H1) to new header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

# Data values for new columns I (I0) and J (IF)
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 9

$ws.Range("I5").Value = 7
$ws.Range("J5").Value = 7

$ws.Range("I6").Value = 8
$ws.Range("J6").Value = 8
